$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 39
$ws.Range("H39").Value = 3854.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 3854.5
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 11563.5
$ws.Range("M39").Value = ""
$ws.Range("N39").Value = -12155.5

# Row 40
$ws.Range("H40").Value = 3799.6
$ws.Range("I40").Value = 3999.5
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 3999.5
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -3824.5
$ws.Range("N40").Value = -3350

# Row 88
$ws.Range("H88").Value = 2846.6667
$ws.Range("I88").Value = 2986.6667
$ws.Range("J88").Value = 2426.6667
$ws.Range("K88").Value = 2986.6667
$ws.Range("L88").Value = 2426.6667
$ws.Range("M88").Value = -2580.6667
$ws.Range("N88").Value = -3238.6667

# Row 91
$ws.Range("H91").Value = 2846.6667
$ws.Range("I91").Value = 2986.6667
$ws.Range("J91").Value = 2426.6667
$ws.Range("K91").Value = 2986.6667
$ws.Range("L91").Value = 2426.6667
$ws.Range("M91").Value = -1582.6667
$ws.Range("N91").Value = -5234.6667

# Row 98
$ws.Range("H98").Value = 9458110
$ws.Range("I98").Value = 11397931
$ws.Range("J98").Value = 1487.5
$ws.Range("K98").Value = 11397931
$ws.Range("L98").Value = 1487.5
$ws.Range("M98").Value = -11396433
$ws.Range("N98").Value = -4483.5

# Row 122
$ws.Range("H122").Value = 9458110
$ws.Range("I122").Value = 11397931
$ws.Range("J122").Value = 1487.5
$ws.Range("K122").Value = 34193793
$ws.Range("L122").Value = 4462.5
$ws.Range("M122").Value = -34191343
$ws.Range("N122").Value = -9362.5

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 3467.5134
$ws.Range("I74").Value = 371.08334
$ws.Range("J74").Value = 9184
$ws.Range("K74").Value = 371.08334
$ws.Range("L74").Value = 9184
$ws.Range("M74").Value = 502.91666
$ws.Range("N74").Value = -10932

# Row 77
$ws.Range("H77").Value = 3467.5134
$ws.Range("I77").Value = 371.08334
$ws.Range("J77").Value = 9184
$ws.Range("K77").Value = 1855.4167
$ws.Range("L77").Value = 45920
$ws.Range("M77").Value = 2512.5833
$ws.Range("N77").Value = -54656

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 4546599.5
$ws.Range("I86").Value = 1038.4445
$ws.Range("J86").Value = 25001624
$ws.Range("K86").Value = 1038.4445
$ws.Range("L86").Value = 25001624
$ws.Range("M86").Value = 84.55549999999994
$ws.Range("N86").Value = -25003870

# Row 89
$ws.Range("H89").Value = 4546599.5
$ws.Range("I89").Value = 1038.4445
$ws.Range("J89").Value = 25001624
$ws.Range("K89").Value = 5192.2225
$ws.Range("L89").Value = 125008120
$ws.Range("M89").Value = 423.7775000000001
$ws.Range("N89").Value = -125019352

# Row 134
$ws.Range("H134").Value = 6546608.5
$ws.Range("I134").Value = 7258010
$ws.Range("J134").Value = 1716.8
$ws.Range("K134").Value = 21774030
$ws.Range("L134").Value = 5150.4
$ws.Range("M134").Value = -21771495
$ws.Range("N134").Value = -10220.4

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9142.941999999999
$ws.Range("I31").Value = 1067.4517
$ws.Range("J31").Value = 21063.904
$ws.Range("K31").Value = 1067.4517
$ws.Range("L31").Value = 21063.904
$ws.Range("M31").Value = -772.4517000000001
$ws.Range("N31").Value = -21653.904

# Row 34
$ws.Range("H34").Value = 9142.941999999999
$ws.Range("I34").Value = 1067.4517
$ws.Range("J34").Value = 21063.904
$ws.Range("K34").Value = 1067.4517
$ws.Range("L34").Value = 21063.904
$ws.Range("M34").Value = -865.4517000000001
$ws.Range("N34").Value = -21467.904

# Row 58
$ws.Range("H58").Value = 3272344
$ws.Range("I58").Value = 3509247.2
$ws.Range("J58").Value = 34666.668
$ws.Range("K58").Value = 3509247.2
$ws.Range("L58").Value = 34666.668
$ws.Range("M58").Value = -3509044.2
$ws.Range("N58").Value = -35072.668

# Row 132
$ws.Range("H132").Value = 13335643
$ws.Range("I132").Value = 27779756
$ws.Range("J132").Value = 2615.2307
$ws.Range("K132").Value = 83339268
$ws.Range("L132").Value = 7845.6921
$ws.Range("M132").Value = -83336738
$ws.Range("N132").Value = -12905.6921

# Row 136
$ws.Range("H136").Value = 3272344
$ws.Range("I136").Value = 3509247.2
$ws.Range("J136").Value = 34666.668
$ws.Range("K136").Value = 10527741.6
$ws.Range("L136").Value = 104000.004
$ws.Range("M136").Value = -10525191.6
$ws.Range("N136").Value = -109100.004

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 910.5
$ws.Range("I5").Value = 826.2857
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 2478.8571
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -2366.8571
$ws.Range("N5").Value = -4724

# Row 132
$ws.Range("H132").Value = 30843.5
$ws.Range("I132").Value = 687
$ws.Range("J132").Value = 61000
$ws.Range("K132").Value = 6183
$ws.Range("L132").Value = 549000
$ws.Range("M132").Value = -3653
$ws.Range("N132").Value = -554060

# Row 135
$ws.Range("H135").Value = 910.5
$ws.Range("I135").Value = 826.2857
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 7436.571300000001
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -4901.571300000001
$ws.Range("N135").Value = -18570

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 3010
$ws.Range("I46").Value = 3010
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 3010
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2854
$ws.Range("N46").Value = ""

# Row 52
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""

# Row 126
$ws.Range("H126").Value = 1041.9584
$ws.Range("I126").Value = 814.0909
$ws.Range("J126").Value = 1234.7693
$ws.Range("K126").Value = 2442.2727
$ws.Range("L126").Value = 3704.3079
$ws.Range("M126").Value = 27.72730000000001
$ws.Range("N126").Value = -8644.3079

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 6166.6665
$ws.Range("I82").Value = 3000
$ws.Range("J82").Value = 6800
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 6800
$ws.Range("M82").Value = -2639
$ws.Range("N82").Value = -7522

# Row 85
$ws.Range("H85").Value = 6166.6665
$ws.Range("I85").Value = 3000
$ws.Range("J85").Value = 6800
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 6800
$ws.Range("M85").Value = -1752
$ws.Range("N85").Value = -9296

# Row 132
$ws.Range("H132").Value = 15388462
$ws.Range("I132").Value = 50003500
$ws.Range("J132").Value = 3999.6667
$ws.Range("K132").Value = 150010500
$ws.Range("L132").Value = 11999.0001
$ws.Range("M132").Value = -150007970
$ws.Range("N132").Value = -17059.0001

# Row 136
$ws.Range("H136").Value = 4923.2705
$ws.Range("I136").Value = 6442.125
$ws.Range("J136").Value = 2119.2307
$ws.Range("K136").Value = 19326.375
$ws.Range("L136").Value = 6357.6921
$ws.Range("M136").Value = -16776.375
$ws.Range("N136").Value = -11457.6921

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 39672212
$ws.Range("I132").Value = 25003994
$ws.Range("J132").Value = 51406784
$ws.Range("K132").Value = 75011982
$ws.Range("L132").Value = 154220352
$ws.Range("M132").Value = -75009452
$ws.Range("N132").Value = -154225412

Write-Output "Applied all market price refresh updates"